$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.149.31'
$ws.Range('E2').Value = '  +4.45%  '
$ws.Range('D3').Value = '2.731.90'
$ws.Range('E3').Value = '  +3.06%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '580.56'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '158.15'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +9.18%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.627'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +4.57%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '2.756.91'
$ws.Range('E9').Value = '  +3.19%  '
$ws.Range('E10').Value = '  +3.11%  '
$ws.Range('E11').Value = '  +2.97%  '
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').Value = '3.224.11'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('E15').Value = '  +3.46%  '
$ws.Range('D16').Value = '64.011.07'
$ws.Range('E16').Value = '  +4.23%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.0000155'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +6.51%  '
$ws.Range('D18').Value = '2.753.39'
$ws.Range('E18').Value = '  +3.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.10'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +4.28%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.95'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '364.11'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +2.79%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '7.02'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.543'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.67%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.993'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '66.94'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +4.53%  '
$ws.Range('E26').Value = '  +5.97%  '
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = '0.0₃0921'
$ws.Range('E29').Value = '  +12.50%  '
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('E31').Value = '  +6.26%  '
$ws.Range('E32').Value = '  +13.52%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '174.06'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '20.67'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.19%  '
$ws.Range('E36').Value = '  +5.97%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.47'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +9.52%  '
$ws.Range('E38').Value = '  +6.75%  '
$ws.Range('E39').Value = '  +10.85%  '
$ws.Range('E40').Value = '  +4.29%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.26'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +18.43%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '338.90'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -1.14%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '39.50'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +3.03%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '21.96'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +7.13%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '22.46'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0604'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.64%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.646'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.91%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0259'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +3.12%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '138.08'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('E50').Value = '  +3.12%  '
$ws.Range('E51').Value = '  +0.02%  '
